# Update the LIDAR flight status (column B, "ESTADO") from "PENDIENTE" to
# "VOLADA" for the rows whose flights have now been flown.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(101, 102, 103, 177, 178, 179, 182, 183, 184, 185, 196, 197, 198, 199, 200, 201, 202, 203)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 2).Value = "VOLADA"
}
